$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 20.32277192010044
$ws.Range("D2").Value = 8.410951872629912
$ws.Range("E2").Value = 13.96188322726679
$ws.Range("F2").Value = 39.06149444435218
$ws.Range("G2").Value = 47.62117084523594
$ws.Range("H2").Value = 18.39160454364103
$ws.Range("J2").Value = 10.56361244083334
$ws.Range("L2").Value = 12.547557784578
$ws.Range("N2").Value = 19.12088269987375
$ws.Range("B3").Value = 20.05608844264176
$ws.Range("D3").Value = 8.336272933928003
$ws.Range("E3").Value = 13.84732200560935
$ws.Range("F3").Value = 38.97729172611429
$ws.Range("G3").Value = 47.12819656069399
$ws.Range("H3").Value = 18.38275860173212
$ws.Range("J3").Value = 10.53253069326606
$ws.Range("L3").Value = 12.34254719339289
$ws.Range("N3").Value = 19.17819226176731
$ws.Range("B4").Value = 19.89459987599795
$ws.Range("D4").Value = 8.289404827876316
$ws.Range("E4").Value = 13.77589216511185
$ws.Range("F4").Value = 38.93795848016762
$ws.Range("G4").Value = 46.841555301432
$ws.Range("H4").Value = 18.38228255020807
$ws.Range("J4").Value = 10.51342870951684
$ws.Range("L4").Value = 12.21706634066093
$ws.Range("N4").Value = 19.21563154905846
$ws.Range("B5").Value = 19.82943400860134
$ws.Range("D5").Value = 8.270056109846111
$ws.Range("E5").Value = 13.74652173096295
$ws.Range("F5").Value = 38.92504302474589
$ws.Range("G5").Value = 46.72890923791465
$ws.Range("H5").Value = 18.38333287801872
$ws.Range("J5").Value = 10.50564177109159
$ws.Range("L5").Value = 12.1660932542377
$ws.Range("N5").Value = 19.23145559546908
$ws.Range("B6").Value = 19.81865420605066
$ws.Range("D6").Value = 8.266828361022885
$ws.Range("E6").Value = 13.74162929596962
$ws.Range("F6").Value = 38.92308647355866
$ws.Range("G6").Value = 46.71045936938808
$ws.Range("H6").Value = 18.38358236917166
$ws.Range("J6").Value = 10.50434864084778
$ws.Range("L6").Value = 12.15764076828648
$ws.Range("N6").Value = 19.23411746161657
$ws.Range("B7").Value = 19.89371832890266
$ws.Range("D7").Value = 8.28914488734034
$ws.Range("E7").Value = 13.77549711139877
$ws.Range("F7").Value = 38.93777169138622
$ws.Range("G7").Value = 46.84001910280563
$ws.Range("H7").Value = 18.38229168019006
$ws.Range("J7").Value = 10.51332370091337
$ws.Range("L7").Value = 12.216378164648
$ws.Range("N7").Value = 19.21584265939863
$ws.Range("B8").Value = 20.23039627844654
$ws.Range("D8").Value = 8.385416578730389
$ws.Range("E8").Value = 13.92261542830583
$ws.Range("F8").Value = 39.02989748841628
$ws.Range("G8").Value = 47.44794241329232
$ws.Range("H8").Value = 18.3875251838652
$ws.Range("J8").Value = 10.55290015520269
$ws.Range("L8").Value = 12.47682470333403
$ws.Range("N8").Value = 19.14017656408703
$ws.Range("B9").Value = 20.90520759882992
$ws.Range("D9").Value = 8.565911073574101
$ws.Range("E9").Value = 14.2019730235665
$ws.Range("F9").Value = 39.30838098474015
$ws.Range("G9").Value = 48.76131212362323
$ws.Range("H9").Value = 18.43714677306584
$ws.Range("J9").Value = 10.63029480051526
$ws.Range("L9").Value = 12.98802842638683
$ws.Range("N9").Value = 19.00960275785175
$ws.Range("B10").Value = 21.405439408953
$ws.Range("D10").Value = 8.693151448902372
$ws.Range("E10").Value = 14.40098424447451
$ws.Range("F10").Value = 39.57197335339674
$ws.Range("G10").Value = 49.79140384113237
$ws.Range("H10").Value = 18.49759658161484
$ws.Range("J10").Value = 10.68692736310918
$ws.Range("L10").Value = 13.36038956413498
$ws.Range("N10").Value = 18.92445347256981
$ws.Range("B11").Value = 21.63304250576651
$ws.Range("D11").Value = 8.749806099132899
$ws.Range("E11").Value = 14.49002842489072
$ws.Range("F11").Value = 39.70448377127713
$ws.Range("G11").Value = 50.27217677839214
$ws.Range("H11").Value = 18.53028644436727
$ws.Range("J11").Value = 10.71261978954192
$ws.Range("L11").Value = 13.52837630023341
$ws.Range("N11").Value = 18.88804348991134
$ws.Range("B12").Value = 21.71916032052535
$ws.Range("D12").Value = 8.771077788836818
$ws.Range("E12").Value = 14.52352175171308
$ws.Range("F12").Value = 39.75645026809722
$ws.Range("G12").Value = 50.4558119080881
$ws.Range("H12").Value = 18.54340830042366
$ws.Range("J12").Value = 10.72233683933018
$ws.Range("L12").Value = 13.59173037539505
$ws.Range("N12").Value = 18.8745892015017
$ws.Range("B13").Value = 21.70061765072277
$ws.Range("D13").Value = 8.76650474196542
$ws.Range("E13").Value = 14.51631859373894
$ws.Range("F13").Value = 39.74517926927993
$ws.Range("G13").Value = 50.4161954632234
$ws.Range("H13").Value = 18.54054929933282
$ws.Range("J13").Value = 10.72024467045052
$ws.Range("L13").Value = 13.57809833888188
$ws.Range("N13").Value = 18.87747201044624
$ws.Range("B14").Value = 21.64012931807587
$ws.Range("D14").Value = 8.751559817333012
$ws.Range("E14").Value = 14.49278853351543
$ws.Range("F14").Value = 39.70872340999222
$ws.Range("G14").Value = 50.287253954427
$ws.Range("H14").Value = 18.53135112199972
$ws.Range("J14").Value = 10.71341947167862
$ws.Range("L14").Value = 13.53359398485311
$ws.Range("N14").Value = 18.88692992125124
$ws.Range("B15").Value = 21.60306702904324
$ws.Range("D15").Value = 8.742381726120579
$ws.Range("E15").Value = 14.47834593914623
$ws.Range("F15").Value = 39.68662517162129
$ws.Range("G15").Value = 50.20847369943995
$ws.Range("H15").Value = 18.52581360045647
$ws.Range("J15").Value = 10.70923720152669
$ws.Range("L15").Value = 13.50629841657439
$ws.Range("N15").Value = 18.89276655381667
$ws.Range("B16").Value = 21.39055929122136
$ws.Range("D16").Value = 8.689423744083202
$ws.Range("E16").Value = 14.39513405004622
$ws.Range("F16").Value = 39.56356504853679
$ws.Range("G16").Value = 49.76021440536562
$ws.Range("H16").Value = 18.49556439735504
$ws.Range("J16").Value = 10.68524671940552
$ws.Range("L16").Value = 13.34937807230661
$ws.Range("N16").Value = 18.92687965240484
$ws.Range("B17").Value = 21.26014883854739
$ws.Range("D17").Value = 8.656617241659713
$ws.Range("E17").Value = 14.34369683276811
$ws.Range("F17").Value = 39.49128284424743
$ws.Range("G17").Value = 49.48821777169488
$ws.Range("H17").Value = 18.47833481196771
$ws.Range("J17").Value = 10.67050995581421
$ws.Range("L17").Value = 13.25271189881111
$ws.Range("N17").Value = 18.94840171303264
$ws.Range("B18").Value = 21.18514894457243
$ws.Range("D18").Value = 8.637632589709687
$ws.Range("E18").Value = 14.31397224290685
$ws.Range("F18").Value = 39.45089577101479
$ws.Range("G18").Value = 49.3329294805474
$ws.Range("H18").Value = 18.46891368083384
$ws.Range("J18").Value = 10.66202729356855
$ws.Range("L18").Value = 13.19698336384936
$ws.Range("N18").Value = 18.96099951298753
$ws.Range("B19").Value = 21.15975914980903
$ws.Range("D19").Value = 8.631185098350686
$ws.Range("E19").Value = 14.30388445600559
$ws.Range("F19").Value = 39.43742608150713
$ws.Range("G19").Value = 49.28055552593253
$ws.Range("H19").Value = 18.46580790327108
$ws.Range("J19").Value = 10.65915416138686
$ws.Range("L19").Value = 13.17809431716351
$ws.Range("N19").Value = 18.96530253323139
$ws.Range("B20").Value = 21.27403091247988
$ws.Range("D20").Value = 8.660121516045304
$ws.Range("E20").Value = 14.34918691265459
$ws.Range("F20").Value = 39.49885464709826
$ws.Range("G20").Value = 49.51705379230044
$ws.Range("H20").Value = 18.4801183569978
$ws.Range("J20").Value = 10.67207939337052
$ws.Range("L20").Value = 13.26301592119493
$ws.Range("N20").Value = 18.94608800653802
$ws.Range("B21").Value = 21.65789872606549
$ws.Range("D21").Value = 8.755954493253691
$ws.Range("E21").Value = 14.49970610411614
$ws.Range("F21").Value = 39.71938307327569
$ws.Range("G21").Value = 50.32508580535134
$ws.Range("H21").Value = 18.53403272020541
$ws.Range("J21").Value = 10.7154245432223
$ws.Range("L21").Value = 13.54667346133359
$ws.Range("N21").Value = 18.88414286294585
$ws.Range("B22").Value = 21.90833312679682
$ws.Range("D22").Value = 8.817522365968369
$ws.Range("E22").Value = 14.5967577361455
$ws.Range("F22").Value = 39.87391743975512
$ws.Range("G22").Value = 50.86229063634139
$ws.Range("H22").Value = 18.57359681743677
$ws.Range("J22").Value = 10.74368152574982
$ws.Range("L22").Value = 13.73052651148821
$ws.Range("N22").Value = 18.84560092615136
$ws.Range("B23").Value = 21.77473724698485
$ws.Range("D23").Value = 8.784761675122935
$ws.Range("E23").Value = 14.54508424923576
$ws.Range("F23").Value = 39.79049627366369
$ws.Range("G23").Value = 50.57479885796228
$ws.Range("H23").Value = 18.55208612799562
$ws.Range("J23").Value = 10.72860747422035
$ws.Range("L23").Value = 13.63255924761635
$ws.Range("N23").Value = 18.8659940160939
$ws.Range("B24").Value = 21.26775489944498
$ws.Range("D24").Value = 8.658537618639395
$ws.Range("E24").Value = 14.34670532209632
$ws.Range("F24").Value = 39.49542779227436
$ws.Range("G24").Value = 49.50401363091305
$ws.Range("H24").Value = 18.47931050727622
$ws.Range("J24").Value = 10.6713698825802
$ws.Range("L24").Value = 13.25835794756006
$ws.Range("N24").Value = 18.94713333396653
$ws.Range("B25").Value = 20.72155213141234
$ws.Range("D25").Value = 8.518003037335426
$ws.Range("E25").Value = 14.12744652260084
$ws.Range("F25").Value = 39.22262309822023
$ws.Range("G25").Value = 48.39390768073553
$ws.Range("H25").Value = 18.41950550369552
$ws.Range("J25").Value = 10.60939050891336
$ws.Range("L25").Value = 12.85005185181761
$ws.Range("N25").Value = 19.04302765841837

Write-Host "Updated loading_percent values for 380 kV case"
